{"js": "// Replace each two-digit-divided-by-one-digit expression with its new value.\n// Every \"old\" string below is unique within the document, so an exact\n// (case-sensitive, whole-match) search safely targets only the intended run.\nconst replacements = [\n  [\"72\u00f77=\", \"48\u00f75=\"],\n  [\"85\u00f76=\", \"98\u00f73=\"],\n  [\"11\u00f73=\", \"66\u00f78=\"],\n  [\"26\u00f75=\", \"24\u00f76=\"],\n  [\"58\u00f78=\", \"79\u00f79=\"],\n  [\"19\u00f72=\", \"42\u00f78=\"],\n  [\"14\u00f73=\", \"36\u00f77=\"],\n  [\"93\u00f74=\", \"25\u00f73=\"],\n  [\"15\u00f77=\", \"84\u00f73=\"],\n  [\"42\u00f75=\", \"23\u00f73=\"],\n  [\"48\u00f77=\", \"70\u00f76=\"],\n  [\"36\u00f79=\", \"11\u00f72=\"],\n  [\"89\u00f76=\", \"56\u00f79=\"],\n  [\"86\u00f74=\", \"32\u00f76=\"],\n  [\"95\u00f74=\", \"61\u00f72=\"],\n  [\"53\u00f72=\", \"54\u00f77=\"],\n  [\"22\u00f73=\", \"60\u00f73=\"],\n  [\"20\u00f74=\", \"43\u00f74=\"],\n  [\"69\u00f72=\", \"60\u00f77=\"],\n  [\"27\u00f73=\", \"14\u00f75=\"],\n  [\"40\u00f77=\", \"83\u00f73=\"],\n  [\"52\u00f77=\", \"76\u00f72=\"],\n  [\"72\u00f75=\", \"69\u00f77=\"],\n  [\"20\u00f76=\", \"14\u00f74=\"],\n  [\"34\u00f74=\", \"91\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-divided-by-one-digit expression with its new value.\n# Every \"old\" string below occurs exactly once in the document, so Find/Replace\n# with MatchWholeWord off but exact text matching safely targets only the\n# intended run each time.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @(\"72\u00f77=\", \"48\u00f75=\"),\n    @(\"85\u00f76=\", \"98\u00f73=\"),\n    @(\"11\u00f73=\", \"66\u00f78=\"),\n    @(\"26\u00f75=\", \"24\u00f76=\"),\n    @(\"58\u00f78=\", \"79\u00f79=\"),\n    @(\"19\u00f72=\", \"42\u00f78=\"),\n    @(\"14\u00f73=\", \"36\u00f77=\"),\n    @(\"93\u00f74=\", \"25\u00f73=\"),\n    @(\"15\u00f77=\", \"84\u00f73=\"),\n    @(\"42\u00f75=\", \"23\u00f73=\"),\n    @(\"48\u00f77=\", \"70\u00f76=\"),\n    @(\"36\u00f79=\", \"11\u00f72=\"),\n    @(\"89\u00f76=\", \"56\u00f79=\"),\n    @(\"86\u00f74=\", \"32\u00f76=\"),\n    @(\"95\u00f74=\", \"61\u00f72=\"),\n    @(\"53\u00f72=\", \"54\u00f77=\"),\n    @(\"22\u00f73=\", \"60\u00f73=\"),\n    @(\"20\u00f74=\", \"43\u00f74=\"),\n    @(\"69\u00f72=\", \"60\u00f77=\"),\n    @(\"27\u00f73=\", \"14\u00f75=\"),\n    @(\"40\u00f77=\", \"83\u00f73=\"),\n    @(\"52\u00f77=\", \"76\u00f72=\"),\n    @(\"72\u00f75=\", \"69\u00f77=\"),\n    @(\"20\u00f76=\", \"14\u00f74=\"),\n    @(\"34\u00f74=\", \"91\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n}\n"}
